# Generate Report for Archive
#
# The localization run moved past handoff, so the recorded status for the
# tracked file flips from "Ready for handoff" to "In Translation" on the
# Overview roll-up (zh-cn/de-de status columns) as well as on each
# per-locale detail sheet. Re-fit the Status column(s) afterwards since
# "In Translation" is shorter than "Ready for handoff".

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Columns("E").AutoFit()
$overview.Columns("F").AutoFit()

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Columns("C").AutoFit()

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Columns("C").AutoFit()
